# Actualizo coef de gini(4T23), VAR IPC ABR24, RES FISCAL ABR24 Y COMEX MAR24
$wb = $excel.ActiveWorkbook

$ws1  = $wb.Worksheets.Item("IPC-Seriemensual")
$ws2  = $wb.Worksheets.Item("IPC-DIC-Div")
$ws3  = $wb.Worksheets.Item("IPC-Interanual")
$ws12 = $wb.Worksheets.Item("Aperturas")

# ---------------------------------------------------------------------------
# IPC-Seriemensual: add a new monthly row (Apr-2024 / serial 45383) below the
# existing data (row 88 -> row 89), inheriting row 88's formatting.
# ---------------------------------------------------------------------------
$ws1.Range("A88:E88").Copy() | Out-Null
$ws1.Range("A89:E89").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(89,1).Value = 45383
$ws1.Cells.Item(89,2).Value = 8.8
$ws1.Cells.Item(89,3).Value = 9.9
$ws1.Cells.Item(89,4).Value = 6.3
$ws1.Cells.Item(89,5).Value = 18.4
$ws1.Range("E89").Select() | Out-Null

# ---------------------------------------------------------------------------
# IPC-DIC-Div: refreshed year-over-year series (rows 2-14, cols B-H)
# ---------------------------------------------------------------------------
$ws2.Cells.Item(2,2).Value = 8.8000000000000007
$ws2.Cells.Item(2,3).Value = 9.1999999999999993
$ws2.Cells.Item(2,4).Value = 8.8000000000000007
$ws2.Cells.Item(2,5).Value = 6.3
$ws2.Cells.Item(2,6).Value = 9.1
$ws2.Cells.Item(2,7).Value = 8.1
$ws2.Cells.Item(2,8).Value = 8.5
$ws2.Cells.Item(3,2).Value = 6
$ws2.Cells.Item(3,3).Value = 5.7
$ws2.Cells.Item(3,4).Value = 6.4
$ws2.Cells.Item(3,5).Value = 4.9000000000000004
$ws2.Cells.Item(3,6).Value = 7.7
$ws2.Cells.Item(3,7).Value = 5.4
$ws2.Cells.Item(3,8).Value = 5.5
$ws2.Cells.Item(4,2).Value = 5.5
$ws2.Cells.Item(4,3).Value = 5.4
$ws2.Cells.Item(4,4).Value = 5.3
$ws2.Cells.Item(4,5).Value = 5.6
$ws2.Cells.Item(4,6).Value = 5.3
$ws2.Cells.Item(4,7).Value = 6.2
$ws2.Cells.Item(4,8).Value = 6.4
$ws2.Cells.Item(5,2).Value = 9.6
$ws2.Cells.Item(5,3).Value = 10.199999999999999
$ws2.Cells.Item(5,4).Value = 9.1999999999999993
$ws2.Cells.Item(5,5).Value = 8
$ws2.Cells.Item(5,6).Value = 10.1
$ws2.Cells.Item(5,7).Value = 10.4
$ws2.Cells.Item(5,8).Value = 6.7
$ws2.Cells.Item(6,2).Value = 35.6
$ws2.Cells.Item(6,3).Value = 40.4
$ws2.Cells.Item(6,4).Value = 33.9
$ws2.Cells.Item(6,5).Value = 12.9
$ws2.Cells.Item(6,6).Value = 31.4
$ws2.Cells.Item(6,7).Value = 30.2
$ws2.Cells.Item(6,8).Value = 39.299999999999997
$ws2.Cells.Item(7,2).Value = 6.5
$ws2.Cells.Item(7,3).Value = 6.4
$ws2.Cells.Item(7,4).Value = 6.7
$ws2.Cells.Item(7,5).Value = 5.6
$ws2.Cells.Item(7,6).Value = 6.2
$ws2.Cells.Item(7,7).Value = 6.1
$ws2.Cells.Item(7,8).Value = 8.1
$ws2.Cells.Item(8,2).Value = 9.1
$ws2.Cells.Item(8,3).Value = 9.5
$ws2.Cells.Item(8,4).Value = 9
$ws2.Cells.Item(8,5).Value = 7.9
$ws2.Cells.Item(8,6).Value = 8.3000000000000007
$ws2.Cells.Item(8,7).Value = 9.1
$ws2.Cells.Item(8,8).Value = 8.6
$ws2.Cells.Item(9,2).Value = 6.3
$ws2.Cells.Item(9,3).Value = 5
$ws2.Cells.Item(9,4).Value = 8.6999999999999993
$ws2.Cells.Item(9,5).Value = 3.8
$ws2.Cells.Item(9,6).Value = 7
$ws2.Cells.Item(9,7).Value = 4.5
$ws2.Cells.Item(9,8).Value = 3.6
$ws2.Cells.Item(10,2).Value = 14.2
$ws2.Cells.Item(10,3).Value = 14.5
$ws2.Cells.Item(10,4).Value = 13.6
$ws2.Cells.Item(10,5).Value = 16.100000000000001
$ws2.Cells.Item(10,6).Value = 14.4
$ws2.Cells.Item(10,7).Value = 12.9
$ws2.Cells.Item(10,8).Value = 13.8
$ws2.Cells.Item(11,2).Value = 7.1
$ws2.Cells.Item(11,3).Value = 6.9
$ws2.Cells.Item(11,4).Value = 6.2
$ws2.Cells.Item(11,5).Value = 6.8
$ws2.Cells.Item(11,6).Value = 11.9
$ws2.Cells.Item(11,7).Value = 8.3000000000000007
$ws2.Cells.Item(11,8).Value = 8
$ws2.Cells.Item(12,2).Value = 8.6
$ws2.Cells.Item(12,3).Value = 8.4
$ws2.Cells.Item(12,4).Value = 8.3000000000000007
$ws2.Cells.Item(12,5).Value = 13.9
$ws2.Cells.Item(12,6).Value = 6.6
$ws2.Cells.Item(12,7).Value = 8.8000000000000007
$ws2.Cells.Item(12,8).Value = 12.3
$ws2.Cells.Item(13,2).Value = 7.3
$ws2.Cells.Item(13,3).Value = 7.8
$ws2.Cells.Item(13,4).Value = 7.2
$ws2.Cells.Item(13,5).Value = 5.4
$ws2.Cells.Item(13,6).Value = 6.3
$ws2.Cells.Item(13,7).Value = 6.5
$ws2.Cells.Item(13,8).Value = 7.1
$ws2.Cells.Item(14,2).Value = 5.7
$ws2.Cells.Item(14,3).Value = 5.2
$ws2.Cells.Item(14,4).Value = 6.1
$ws2.Cells.Item(14,5).Value = 5.5
$ws2.Cells.Item(14,6).Value = 4.3
$ws2.Cells.Item(14,7).Value = 7.5
$ws2.Cells.Item(14,8).Value = 7.3

# ---------------------------------------------------------------------------
# IPC-Interanual: refreshed index-level series (rows 2-14, cols B-H)
# ---------------------------------------------------------------------------
$ws3.Cells.Item(2,2).Value = 289.39999999999998
$ws3.Cells.Item(2,3).Value = 292.2
$ws3.Cells.Item(2,4).Value = 288.89999999999998
$ws3.Cells.Item(2,5).Value = 278.89999999999998
$ws3.Cells.Item(2,6).Value = 282.2
$ws3.Cells.Item(2,7).Value = 284.10000000000002
$ws3.Cells.Item(2,8).Value = 293.5
$ws3.Cells.Item(3,2).Value = 293
$ws3.Cells.Item(3,3).Value = 297.5
$ws3.Cells.Item(3,4).Value = 289
$ws3.Cells.Item(3,5).Value = 270.7
$ws3.Cells.Item(3,6).Value = 289.60000000000002
$ws3.Cells.Item(3,7).Value = 290.7
$ws3.Cells.Item(3,8).Value = 310.3
$ws3.Cells.Item(4,2).Value = 272.3
$ws3.Cells.Item(4,3).Value = 268.8
$ws3.Cells.Item(4,4).Value = 278
$ws3.Cells.Item(4,5).Value = 263.89999999999998
$ws3.Cells.Item(4,6).Value = 265.2
$ws3.Cells.Item(4,7).Value = 274
$ws3.Cells.Item(4,8).Value = 282.89999999999998
$ws3.Cells.Item(5,2).Value = 205
$ws3.Cells.Item(5,3).Value = 205.4
$ws3.Cells.Item(5,4).Value = 207.8
$ws3.Cells.Item(5,5).Value = 202.6
$ws3.Cells.Item(5,6).Value = 190.3
$ws3.Cells.Item(5,7).Value = 197.2
$ws3.Cells.Item(5,8).Value = 213
$ws3.Cells.Item(6,2).Value = 311.60000000000002
$ws3.Cells.Item(6,3).Value = 321.3
$ws3.Cells.Item(6,4).Value = 296.7
$ws3.Cells.Item(6,5).Value = 279
$ws3.Cells.Item(6,6).Value = 326
$ws3.Cells.Item(6,7).Value = 317.10000000000002
$ws3.Cells.Item(6,8).Value = 318.60000000000002
$ws3.Cells.Item(7,2).Value = 293.39999999999998
$ws3.Cells.Item(7,3).Value = 286.5
$ws3.Cells.Item(7,4).Value = 299.10000000000002
$ws3.Cells.Item(7,5).Value = 288.7
$ws3.Cells.Item(7,6).Value = 296.3
$ws3.Cells.Item(7,7).Value = 293.7
$ws3.Cells.Item(7,8).Value = 319.2
$ws3.Cells.Item(8,2).Value = 341.1
$ws3.Cells.Item(8,3).Value = 340.3
$ws3.Cells.Item(8,4).Value = 340.8
$ws3.Cells.Item(8,5).Value = 341.9
$ws3.Cells.Item(8,6).Value = 340.1
$ws3.Cells.Item(8,7).Value = 354
$ws3.Cells.Item(8,8).Value = 335.7
$ws3.Cells.Item(9,2).Value = 325.89999999999998
$ws3.Cells.Item(9,3).Value = 331.2
$ws3.Cells.Item(9,4).Value = 326.60000000000002
$ws3.Cells.Item(9,5).Value = 333.1
$ws3.Cells.Item(9,6).Value = 332.7
$ws3.Cells.Item(9,7).Value = 291.60000000000002
$ws3.Cells.Item(9,8).Value = 299.5
$ws3.Cells.Item(10,2).Value = 369.1
$ws3.Cells.Item(10,3).Value = 377.1
$ws3.Cells.Item(10,4).Value = 376.8
$ws3.Cells.Item(10,5).Value = 377.3
$ws3.Cells.Item(10,6).Value = 341.2
$ws3.Cells.Item(10,7).Value = 324.2
$ws3.Cells.Item(10,8).Value = 318.10000000000002
$ws3.Cells.Item(11,2).Value = 283.60000000000002
$ws3.Cells.Item(11,3).Value = 295.2
$ws3.Cells.Item(11,4).Value = 273.89999999999998
$ws3.Cells.Item(11,5).Value = 288.7
$ws3.Cells.Item(11,6).Value = 261.3
$ws3.Cells.Item(11,7).Value = 276
$ws3.Cells.Item(11,8).Value = 282.60000000000002
$ws3.Cells.Item(12,2).Value = 214.5
$ws3.Cells.Item(12,3).Value = 236.2
$ws3.Cells.Item(12,4).Value = 190.3
$ws3.Cells.Item(12,5).Value = 189.7
$ws3.Cells.Item(12,6).Value = 218.5
$ws3.Cells.Item(12,7).Value = 190
$ws3.Cells.Item(12,8).Value = 180.9
$ws3.Cells.Item(13,2).Value = 263.60000000000002
$ws3.Cells.Item(13,3).Value = 262.8
$ws3.Cells.Item(13,4).Value = 271.2
$ws3.Cells.Item(13,5).Value = 248.3
$ws3.Cells.Item(13,6).Value = 250.3
$ws3.Cells.Item(13,7).Value = 251.3
$ws3.Cells.Item(13,8).Value = 268.10000000000002
$ws3.Cells.Item(14,2).Value = 360.2
$ws3.Cells.Item(14,3).Value = 340.4
$ws3.Cells.Item(14,4).Value = 367
$ws3.Cells.Item(14,5).Value = 378
$ws3.Cells.Item(14,6).Value = 386.3
$ws3.Cells.Item(14,7).Value = 387.5
$ws3.Cells.Item(14,8).Value = 418.5

# ---------------------------------------------------------------------------
# Aperturas (coeficiente de Gini por regiones, 4T23): rows 2-10, cols B-H
# ---------------------------------------------------------------------------
$ws12.Cells.Item(2,2).Value = 5.2
$ws12.Cells.Item(2,3).Value = 4.5
$ws12.Cells.Item(2,4).Value = 2.9
$ws12.Cells.Item(2,5).Value = 4.4000000000000004
$ws12.Cells.Item(2,6).Value = 3.3
$ws12.Cells.Item(2,7).Value = 5.0999999999999996
$ws12.Cells.Item(2,8).Value = 4.702169367571174
$ws12.Cells.Item(3,2).Value = 4.9000000000000004
$ws12.Cells.Item(3,3).Value = 5.7
$ws12.Cells.Item(3,4).Value = 3.1
$ws12.Cells.Item(3,5).Value = 7.6
$ws12.Cells.Item(3,6).Value = 4.5999999999999996
$ws12.Cells.Item(3,7).Value = 2.8
$ws12.Cells.Item(3,8).Value = 5.1412244225463866
$ws12.Cells.Item(4,2).Value = 8.4
$ws12.Cells.Item(4,3).Value = 10
$ws12.Cells.Item(4,4).Value = 9.5
$ws12.Cells.Item(4,5).Value = 11.6
$ws12.Cells.Item(4,6).Value = 9.6
$ws12.Cells.Item(4,7).Value = 11.4
$ws12.Cells.Item(4,8).Value = 9.3770117378148221
$ws12.Cells.Item(5,2).Value = 4.0999999999999996
$ws12.Cells.Item(5,3).Value = 4.7
$ws12.Cells.Item(5,4).Value = 1.2
$ws12.Cells.Item(5,5).Value = 4.3
$ws12.Cells.Item(5,6).Value = 7.5
$ws12.Cells.Item(5,7).Value = 4.8
$ws12.Cells.Item(5,8).Value = 4.3410982481665528
$ws12.Cells.Item(6,2).Value = -2.9
$ws12.Cells.Item(6,3).Value = 0.4
$ws12.Cells.Item(6,4).Value = 1.7
$ws12.Cells.Item(6,5).Value = -4.9000000000000004
$ws12.Cells.Item(6,6).Value = -7.8
$ws12.Cells.Item(6,7).Value = 1.1000000000000001
$ws12.Cells.Item(6,8).Value = -1.7959996380280674
$ws12.Cells.Item(7,2).Value = 14.3
$ws12.Cells.Item(7,3).Value = 16.899999999999999
$ws12.Cells.Item(7,4).Value = 14.7
$ws12.Cells.Item(7,5).Value = 24.4
$ws12.Cells.Item(7,6).Value = 18.5
$ws12.Cells.Item(7,7).Value = 12.9
$ws12.Cells.Item(7,8).Value = 15.955457669427076
$ws12.Cells.Item(8,2).Value = 4.2
$ws12.Cells.Item(8,3).Value = 1.7
$ws12.Cells.Item(8,4).Value = 1.8
$ws12.Cells.Item(8,5).Value = -0.2
$ws12.Cells.Item(8,6).Value = -0.5
$ws12.Cells.Item(8,7).Value = 3.5
$ws12.Cells.Item(8,8).Value = 2.5906587853438356
$ws12.Cells.Item(9,2).Value = -0.3
$ws12.Cells.Item(9,3).Value = 3.9
$ws12.Cells.Item(9,4).Value = 4.3
$ws12.Cells.Item(9,5).Value = 3.9
$ws12.Cells.Item(9,6).Value = 4.2
$ws12.Cells.Item(9,7).Value = 6.9
$ws12.Cells.Item(9,8).Value = 2.1793211766671039
$ws12.Cells.Item(10,2).Value = 4.9000000000000004
$ws12.Cells.Item(10,3).Value = 3.3
$ws12.Cells.Item(10,4).Value = 2.6
$ws12.Cells.Item(10,5).Value = 4.2
$ws12.Cells.Item(10,6).Value = 1.5
$ws12.Cells.Item(10,7).Value = 3.2
$ws12.Cells.Item(10,8).Value = 3.8826647884975696

# ---------------------------------------------------------------------------
# Active-tab bookkeeping: the workbook was last left on "Aperturas"
# (previously it was "EMAE") -- activating it flips tabSelected + activeTab.
# ---------------------------------------------------------------------------
$ws12.Activate()

Write-Host "done"
